$wb = $excel.ActiveWorkbook

$journal = $wb.Worksheets.Item("Journal")
$testData = $wb.Worksheets.Item("PotenotTask test data")

# --- Fix typo "PotenotTaskServece" -> "PotenotTaskService" everywhere it is used ---
$journal.Range("D4").Value = "PotenotTaskService"
$journal.Range("D6").Value = "PotenotTaskService"
$journal.Range("D7").Value = "PotenotTaskService"
$journal.Range("D8").Value = "PotenotTaskService"
$journal.Range("D9").Value = "PotenotTaskService"

# --- Add two new journal entries for the service-layer test work ---
# Row 10: tests for PotenotTaskMapperImpl, 2024-09-16, 45 minutes
$journal.Range("A10").Value = "Разработка тестов для PotenotTaskMapperImpl"
$journal.Range("B10").Value = 45551
$journal.Range("C10").Value = 0.03125
$journal.Range("D10").Value = "PotenotTaskService"

# Row 11: tests for the service layer, 2024-09-16, 1 hour
$journal.Range("A11").Value = "Разработка тестов для слоя servise"
$journal.Range("B11").Value = 45551
$journal.Range("C11").Value = 0.0416666666666667
$journal.Range("D11").Value = "PotenotTaskService"

# --- New computed totals columns (AE/AF) on the test-data sheet ---
# Row 2, 4 and 5 pick up the neighbouring bordered cell style, row 3 stays unstyled
$testData.Range("AC2").Copy($testData.Range("AE2"))
$testData.Range("AC2").Copy($testData.Range("AF2"))
$testData.Range("AE2").Value = 137114
$testData.Range("AF2").Value = 209238

$testData.Range("AE3").Value = 18515
$testData.Range("AF3").Value = 143314

$testData.Range("AC4").Copy($testData.Range("AE4"))
$testData.Range("AC4").Copy($testData.Range("AF4"))
$testData.Range("AE4").Value = -709005
$testData.Range("AF4").Value = 846436

$testData.Range("AC5").Copy($testData.Range("AE5"))
$testData.Range("AC5").Copy($testData.Range("AF5"))
$testData.Range("AE5").Value = 150295
$testData.Range("AF5").Value = -176044

# --- View/selection state ---
$testData.Activate()
$testData.Range("M5").Select()

$journal.Activate()
$journal.Range("C12").Select()
